$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell B1 from "a" to "A"
$ws.Range("B1").Value = "A"

# Move the active selection to B2
$ws.Range("B2").Select()
